# Apply the "10 more IMG annotations of genomes" edit: fill in the
# newly-annotated columns (J/K/L/M/N/P/Q/S/T) for rows 70-80.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J70").Value = "chitobiose, glucose, fructose, galactose, mannose, sucrose, galacturonate, maltose"
$ws.Range("L70").Value = "carotenoid synthesis"
$ws.Range("N70").Value = "ammonia_assimilation"
$ws.Range("P70").Value = "sulfate_red_ass"
$ws.Range("Q70").Value = "heme, LPS export, lipoprotein release, oligopepteide, phosphate, phospholipid, phosphonate, spermidine/putrescine"
$ws.Range("S70").Value = "flagellum with type III secretion"
$ws.Range("T70").Value = "Oxidative phosphorylation"

$ws.Range("J71").Value = "chitobiose, glycolate, cellobiose, glucose"
$ws.Range("L71").Value = "some carotenoid synthesis"
$ws.Range("N71").Value = "ammonia_assimilation"
$ws.Range("Q71").Value = "heme, LPS export, lipoprotein release, phospholipid/cholesterol"
$ws.Range("T71").Value = "Oxidative phosphorylation"

$ws.Range("K72").Value = "glycolate, cellobiose"
$ws.Range("L72").Value = "two carotenoid genes"
$ws.Range("N72").Value = "ammonia_assimilation"
$ws.Range("Q72").Value = "heme, iron, LPS export, lipoprotein release, peptide/nickel, phospholipid/cholesterol"
$ws.Range("T72").Value = "Oxidative phosphorylation"

$ws.Range("J73").Value = "carbon fixation (Rubisco), chitobiose, glucose, glycolate, starch/glycogen, maltose"
$ws.Range("L73").Value = "some carotenoid synthesis"
$ws.Range("N73").Value = "nitrate_red_ass, nitrogen fixation, ammonia_assimilation"
$ws.Range("Q73").Value = "methionine, arginine/lysine/histidine/glutamine, bicarbonate, biotin, branched amino, carbohydrate, cobalt/nickel, L-amino acid, iron, LPS export, LPS transport, molybdate, neutral amino, nitrate/nitrite, phosphate, phospholipid/cholesterol, phosphonate, chitobiose, spermidine/putrescine, sulfate, urea"
$ws.Range("T73").Value = "Oxidative phosphorylation, photosynthesis"

$ws.Range("J74").Value = "two carotenoid genes"
$ws.Range("L74").Value = "glucose, galactose, cellulose degradation, maltose, starch/glycogen"
$ws.Range("N74").Value = "ammonia_assimilation"
$ws.Range("P74").Value = "sulfate_red_ass, alkanesulfonate"
$ws.Range("Q74").Value = "heme, LPS export, lipoprotein, molybdate, oligopeptide, phospholipid/cholesterol"
$ws.Range("T74").Value = "Oxidative phosphorylation"

$ws.Range("J75").Value = "Some carotenoid genes"
$ws.Range("N75").Value = "ammonia_assimilation"
$ws.Range("Q75").Value = "iron, LPS export, LPS transport, lipoprotein release, peptide/nickel, phospholipid/cholesterol"
$ws.Range("T75").Value = "Oxidative phosphorylation"

$ws.Range("J76").Value = "two carotenoid genes"
$ws.Range("L76").Value = "cellulose degradation"
$ws.Range("N76").Value = "ammonia_assimilation"
$ws.Range("P76").Value = "one sulfate_red_ass gene"
$ws.Range("Q76").Value = "LPS export, lipoprotein release, phosphate, phospholipid/cholesterol"
$ws.Range("T76").Value = "Oxidative phosphorylation"

$ws.Range("J77").Value = "two carotenoid genes"
$ws.Range("L77").Value = "chitin degradation, glycolate, cellulose degradation"
$ws.Range("N77").Value = "ammonia_assimilation"
$ws.Range("Q77").Value = "polysaccharide phosphate, heme, LPS export, LPS transport, lipoprotein release, phosphate, phospholipid/cholesterol"
$ws.Range("T77").Value = "Oxidative phosphorylation"

$ws.Range("J78").Value = "three carotenoid genes"
$ws.Range("L78").Value = "chitobiose, glucose, galactose, fructose, starch/glycogen, cellulose degradation, maltose"
$ws.Range("N78").Value = "ammonia_assimilation"
$ws.Range("P78").Value = "alkanesulfonate"
$ws.Range("Q78").Value = "heme, iron (III), LPS export, lipoprotein release, phosphate, phospholipid/cholesterol"
$ws.Range("S78").Value = "one chemotaxis protein"
$ws.Range("T78").Value = "Oxidative phosphorylation, luciferase gene?"

$ws.Range("J79").Value = "two carotenoid genes"
$ws.Range("L79").Value = "galactose, melibiose, cellobiose"
$ws.Range("Q79").Value = "iron, LPS export, peptide/nickel, phospholipid/cholesterol"
$ws.Range("T79").Value = "Oxidative phosphorylation"

$ws.Range("J80").Value = "two carotenoid genes"
$ws.Range("L80").Value = "glucose, galacturonate, cellobiose, starch/glycogen, maltose"
$ws.Range("N80").Value = "ammonia_assimilation, one nitritre reductase"
$ws.Range("P80").Value = "sulfate_red_ass"
$ws.Range("Q80").Value = "branched amino, carbohydrate, LPS export, oligopeptide, peptide/nickel, phospholipid/cholesterol"
$ws.Range("S80").Value = "flagellum with type III secretion"
$ws.Range("T80").Value = "Oxidative phosphorylation"

# Match the scrolled/selected viewport recorded in the saved worksheet
# (the last cell touched while annotating was L80).
$ws.Range("L80").Select()
$excel.ActiveWindow.ScrollRow = 62
$excel.ActiveWindow.ScrollColumn = 1
